# Auto-generated Excel COM-interop script to apply numeric updates
# to the Sagittarius_Profits workbook sheets, per the provided diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 737.3333
$ws.Range("I33").Value = 807.0769
$ws.Range("K33").Value = 807.0769
$ws.Range("M33").Value = -578.0769
$ws.Range("H74").Value = 91327.95
$ws.Range("I74").Value = 104472.625
$ws.Range("K74").Value = 104472.625
$ws.Range("M74").Value = -103536.625
$ws.Range("H77").Value = 91327.95
$ws.Range("I77").Value = 104472.625
$ws.Range("K77").Value = 522363.125
$ws.Range("M77").Value = -517683.125
$ws.Range("H112").Value = 1192.5264
$ws.Range("J112").Value = 1175.4445
$ws.Range("L112").Value = 3526.3335
$ws.Range("N112").Value = -5742.333500000001
$ws.Range("H116").Value = 4500
$ws.Range("I116").Value = 4500
$ws.Range("K116").Value = 4500
$ws.Range("M116").Value = -1058
$ws.Range("H132").Value = 1128.091
$ws.Range("I132").Value = 1052.4
$ws.Range("K132").Value = 3157.2
$ws.Range("M132").Value = -627.2000000000003

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 4814
$ws.Range("I28").Value = 4814
$ws.Range("K28").Value = 4814
$ws.Range("M28").Value = -4622
$ws.Range("H45").Value = 1857.8572
$ws.Range("I45").Value = 1848.4615
$ws.Range("K45").Value = 1848.4615
$ws.Range("M45").Value = -1471.4615
$ws.Range("H53").Value = 19939
$ws.Range("I53").Value = 19939
$ws.Range("K53").Value = 19939
$ws.Range("M53").Value = -19257
$ws.Range("H74").Value = 2292.1724
$ws.Range("I74").Value = 2041.8
$ws.Range("J74").Value = 3857
$ws.Range("K74").Value = 2041.8
$ws.Range("L74").Value = 3857
$ws.Range("M74").Value = -1167.8
$ws.Range("N74").Value = -5605
$ws.Range("H77").Value = 2292.1724
$ws.Range("I77").Value = 2041.8
$ws.Range("J77").Value = 3857
$ws.Range("K77").Value = 10209
$ws.Range("L77").Value = 19285
$ws.Range("M77").Value = -5841
$ws.Range("N77").Value = -28021
$ws.Range("H99").Value = 4814
$ws.Range("I99").Value = 4814
$ws.Range("K99").Value = 4814
$ws.Range("M99").Value = -1819
$ws.Range("H110").Value = 246.25
$ws.Range("I110").Value = 246.25
$ws.Range("K110").Value = 246.25
$ws.Range("M110").Value = 1798.75
$ws.Range("H122").Value = 1643.9286
$ws.Range("J122").Value = 2679.8
$ws.Range("L122").Value = 8039.400000000001
$ws.Range("N122").Value = -12939.4
$ws.Range("H128").Value = 74999
$ws.Range("J128").Value = 74999
$ws.Range("L128").Value = 74999
$ws.Range("N128").Value = -84959

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2331.4707
$ws.Range("I20").Value = 2417.8
$ws.Range("J20").Value = 1684
$ws.Range("K20").Value = 2417.8
$ws.Range("L20").Value = 1684
$ws.Range("M20").Value = -2170.8
$ws.Range("N20").Value = -2178
$ws.Range("H94").Value = 512.1111
$ws.Range("I94").Value = 526.125
$ws.Range("K94").Value = 526.125
$ws.Range("M94").Value = -75.125
$ws.Range("H99").Value = 3564.6667
$ws.Range("I99").Value = 2940
$ws.Range("K99").Value = 2940
$ws.Range("M99").Value = -1442

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 9468.259
$ws.Range("I22").Value = 207.61905
$ws.Range("K22").Value = 207.61905
$ws.Range("M22").Value = 142.38095
$ws.Range("H81").Value = 100000
$ws.Range("J81").Value = 100000
$ws.Range("L81").Value = 100000
$ws.Range("N81").Value = -101996
$ws.Range("H84").Value = 100000
$ws.Range("J84").Value = 100000
$ws.Range("L84").Value = 300000
$ws.Range("N84").Value = -309984
$ws.Range("H134").Value = 2035.1471
$ws.Range("I134").Value = 1793.3
$ws.Range("K134").Value = 5379.9
$ws.Range("M134").Value = -2844.9

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 500
$ws.Range("I48").Value = 500
$ws.Range("K48").Value = 1500
$ws.Range("M48").Value = -1250
$ws.Range("H59").Value = 1174.5
$ws.Range("I59").Value = 1333
$ws.Range("J59").Value = 699
$ws.Range("K59").Value = 3999
$ws.Range("L59").Value = 2097
$ws.Range("M59").Value = -3459
$ws.Range("N59").Value = -3177
$ws.Range("H113").Value = 1559.8948
$ws.Range("J113").Value = 1423.3334
$ws.Range("L113").Value = 4270.0002
$ws.Range("N113").Value = -8610.0002
$ws.Range("H132").Value = 2864.95
$ws.Range("I132").Value = 2376.4119
$ws.Range("J132").Value = 5633.3335
$ws.Range("K132").Value = 21387.7071
$ws.Range("L132").Value = 50700.0015
$ws.Range("M132").Value = -18857.7071
$ws.Range("N132").Value = -55760.0015
$ws.Range("H137").Value = 2874.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 48832.332
$ws.Range("I5").Value = 48832.332
$ws.Range("K5").Value = 48832.332
$ws.Range("M5").Value = -48720.332
$ws.Range("H58").Value = 123265
$ws.Range("J58").Value = 9897
$ws.Range("L58").Value = 9897
$ws.Range("N58").Value = -10451
$ws.Range("H70").Value = 5342.9287
$ws.Range("I70").Value = 4890.7144
$ws.Range("K70").Value = 4890.7144
$ws.Range("M70").Value = -4620.7144
$ws.Range("H73").Value = 5342.9287
$ws.Range("I73").Value = 4890.7144
$ws.Range("K73").Value = 4890.7144
$ws.Range("M73").Value = -3954.7144
$ws.Range("H97").Value = 988.7778
$ws.Range("I97").Value = 848.8333
$ws.Range("J97").Value = 1268.6666
$ws.Range("K97").Value = 848.8333
$ws.Range("L97").Value = 1268.6666
$ws.Range("M97").Value = -352.8333
$ws.Range("N97").Value = -2260.6666

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2612500
$ws.Range("I2").Value = 200000
$ws.Range("J2").Value = 5025000
$ws.Range("K2").Value = 200000
$ws.Range("L2").Value = 5025000
$ws.Range("M2").Value = -199888
$ws.Range("N2").Value = -5025224
$ws.Range("H7").Value = 8285.608
$ws.Range("I7").Value = 6936
$ws.Range("K7").Value = 6936
$ws.Range("M7").Value = -6824
$ws.Range("H22").Value = 3208.125
$ws.Range("J22").Value = 2657.8
$ws.Range("L22").Value = 2657.8
$ws.Range("N22").Value = -3247.8
$ws.Range("H27").Value = 3208.125
$ws.Range("J27").Value = 2657.8
$ws.Range("L27").Value = 2657.8
$ws.Range("N27").Value = -2871.8
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("N57").ClearContents()
$ws.Range("H61").Value = 8042.421
$ws.Range("I61").Value = 6187.6
$ws.Range("K61").Value = 6187.6
$ws.Range("M61").Value = -5985.6
$ws.Range("H68").Value = 2938.111
$ws.Range("I68").Value = 2300.1
$ws.Range("J68").Value = 3735.625
$ws.Range("K68").Value = 2300.1
$ws.Range("L68").Value = 3735.625
$ws.Range("M68").Value = -1551.1
$ws.Range("N68").Value = -5233.625
$ws.Range("H71").Value = 2938.111
$ws.Range("I71").Value = 2300.1
$ws.Range("J71").Value = 3735.625
$ws.Range("K71").Value = 11500.5
$ws.Range("L71").Value = 18678.125
$ws.Range("M71").Value = -7756.5
$ws.Range("N71").Value = -26166.125
$ws.Range("H113").Value = 8042.421
$ws.Range("I113").Value = 6187.6
$ws.Range("K113").Value = 6187.6
$ws.Range("M113").Value = -4017.6
$ws.Range("H122").Value = 6053.5356
$ws.Range("I122").Value = 5505.421
$ws.Range("J122").Value = 7210.6665
$ws.Range("K122").Value = 16516.263
$ws.Range("L122").Value = 21631.9995
$ws.Range("M122").Value = -14066.263
$ws.Range("N122").Value = -26531.9995
$ws.Range("H126").Value = 8285.608
$ws.Range("I126").Value = 6936
$ws.Range("K126").Value = 20808
$ws.Range("M126").Value = -18338
$ws.Range("H136").Value = 10753.143
$ws.Range("I136").Value = 10753.143
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 32259.429
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -29709.429
$ws.Range("N136").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 10999.5
$ws.Range("J74").Value = 10999.5
$ws.Range("L74").Value = 10999.5
$ws.Range("N74").Value = -12871.5
$ws.Range("H77").Value = 10999.5
$ws.Range("J77").Value = 10999.5
$ws.Range("L77").Value = 32998.5
$ws.Range("N77").Value = -42358.5
$ws.Range("H81").Value = 7046.1665
$ws.Range("I81").Value = 4498.4
$ws.Range("K81").Value = 8996.799999999999
$ws.Range("M81").Value = -7935.799999999999
$ws.Range("H84").Value = 7046.1665
$ws.Range("I84").Value = 4498.4
$ws.Range("K84").Value = 44984
$ws.Range("M84").Value = -39680
$ws.Range("H124").Value = 43425.285
$ws.Range("J124").Value = 43425.285
$ws.Range("L124").Value = 43425.285
$ws.Range("N124").Value = -53245.285
$ws.Range("H136").Value = 1742.1052
$ws.Range("I136").Value = 1742.1052
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5226.3156
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2676.3156
$ws.Range("N136").ClearContents()

Write-Output "Applied 232 value updates and 4 clears."